# Weekly update: a new day's price observation is inserted at row 203
# (pushing all subsequent rows for this product down by one row), and the
# workbook's used range grows from A1:R323 to A1:R324.
#
# The newly inserted row 203 keeps the same Mercado/Region/Categoria/etc.
# values as before, but carries a new date (Fecha) and a new Volumen value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 203, shifting rows 203:323 down to 204:324.
$ws.Rows("203:203").Insert()

# Populate the new row 203 with the latest weekly observation.
$ws.Range("A203").Value = 3
$ws.Range("B203").Value = "Femacal de La Calera"
$ws.Range("C203").Value = "Coquimbo"
$ws.Range("D203").Value = 44719
$ws.Range("E203").Value = 5
$ws.Range("F203").Value = 100112039
$ws.Range("G203").Value = "Ciboulette"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 120
$ws.Range("K203").Value = 1500
$ws.Range("L203").Value = 1500
$ws.Range("M203").Value = 1500
$ws.Range("N203").Value = '$/docena de atados'
$ws.Range("O203").Value = "Provincia de Quillota"
$ws.Range("P203").Value = 500
$ws.Range("Q203").Value = 3
$ws.Range("R203").Value = "Hortaliza"
